$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.730.28"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "2.476.53"
$ws.Range("E3").Value = "  +0.37%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.73"
$ws.Range("E5").Value = "  +1.22%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "93.18"
$ws.Range("E6").Value = "  +1.20%  "

$ws.Range("E7").Value = "  +0.67%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("E10").Value = "  +10.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "33.25"
$ws.Range("E11").Value = "  +2.95%  "

$ws.Range("E12").Value = "  +0.57%  "

$ws.Range("D13").Value = "2.859.72"
$ws.Range("E13").Value = "  +0.40%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.90"
$ws.Range("E14").Value = "  +0.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.77"
$ws.Range("E15").Value = "  -0.99%  "

$ws.Range("D16").Value = "2.483.07"
$ws.Range("E16").Value = "  +0.75%  "

$ws.Range("E17").Value = "  +2.43%  "

$ws.Range("D18").Value = "41.717.65"
$ws.Range("E18").Value = "  +0.48%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.48"
$ws.Range("E19").Value = "  +0.25%  "

$ws.Range("E20").Value = "  +1.12%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.23"
$ws.Range("E21").Value = "  +0.15%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.31"
$ws.Range("E22").Value = "  +1.79%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "239.53"
$ws.Range("E23").Value = "  +1.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.75"
$ws.Range("E24").Value = "  +1.05%  "

$ws.Range("E25").Value = "  +2.70%  "

$ws.Range("E26").Value = "  +0.03%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.71"

$ws.Range("E28").Value = "  +1.35%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.82"
$ws.Range("E29").Value = "  +1.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "36.17"
$ws.Range("E30").Value = "  +1.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.50"
$ws.Range("E31").Value = "  +1.76%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.52"
$ws.Range("E32").Value = "  +1.51%  "

$ws.Range("E33").Value = "  -0.05%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0767"
$ws.Range("E35").Value = "  +1.16%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.48"
$ws.Range("E36").Value = "  +0.59%  "

$ws.Range("E37").Value = "  +5.75%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.93"
$ws.Range("E38").Value = "  +2.16%  "

$ws.Range("E39").Value = "  +1.57%  "

$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.05"
$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("E42").Value = "  +10.88%  "

$ws.Range("D43").Value = "1.993.09"

$ws.Range("E44").Value = "  +0.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.77"
$ws.Range("E45").Value = "  +1.62%  "

$ws.Range("E46").Value = "  +2.19%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.47"
$ws.Range("E47").Value = "  +4.72%  "

$ws.Range("D48").Value = "2.716.87"
$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "97.60"
$ws.Range("E49").Value = "  +0.74%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.15"
$ws.Range("E50").Value = "  +3.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "67.21"
$ws.Range("E51").Value = "  +0.11%  "
